$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.95991976981970684
$ws.Range("AA1").Value = 0.98085046059712833
$ws.Range("AH1").Value = 0.83147037874169105
$ws.Range("BE1").Value = 0.80527036510961825
$ws.Range("BO1").Value = 0.93839421162811698
$ws.Range("AE2").Value = 0.7506567917394763
$ws.Range("AV2").Value = 0.90304343121070085
$ws.Range("BK2").Value = 0.84379677820922372
$ws.Range("O3").Value = 0.99656285882594542
$ws.Range("AU3").Value = 0.77022442182313344
$ws.Range("T4").Value = 0.65440092083640122
$ws.Range("V4").Value = 0.92244748101769891
$ws.Range("AF4").Value = 0.98462814948436594
$ws.Range("AE5").Value = 0.7092265982583339
$ws.Range("BE5").Value = 0.96653699472827892
$ws.Range("BP5").Value = 0.88588523382743689
$ws.Range("AC6").Value = 0.8974640501488701
$ws.Range("AY6").Value = 0.96867377715060821
$ws.Range("AZ6").Value = 0.90587327298011355
$ws.Range("AA7").Value = 0.96483776852120118
$ws.Range("AW7").Value = 0.96613776097583814
$ws.Range("BK7").Value = 0.87265710627430115
$ws.Range("BP7").Value = 0.8720803951692192
$ws.Range("AH8").Value = 0.53905706255385932
$ws.Range("H9").Value = 0.9266718365048181
$ws.Range("V9").Value = 0.87686217600983807
$ws.Range("AD9").Value = 0.95143292880983554
$ws.Range("H10").Value = 0.90384036773793897
$ws.Range("S10").Value = 0.74365967189862359
$ws.Range("J12").Value = 0.99682583520566492
$ws.Range("N12").Value = 0.9073373489074037
$ws.Range("I13").Value = 0.80933802812260014
$ws.Range("AG13").Value = 0.79690810266863021
$ws.Range("AD15").Value = 0.82553508876880377
$ws.Range("AM15").Value = 0.96252250563021879
$ws.Range("BJ15").Value = 0.59085243308700441
$ws.Range("N16").Value = 0.82069233568954347
$ws.Range("BI16").Value = 0.89163168649589497
$ws.Range("AL17").Value = 0.98105192241968631
$ws.Range("BK17").Value = 0.94120616421081138
$ws.Range("AF18").Value = 0.535403455510715
$ws.Range("T19").Value = 0.8673822094405943
$ws.Range("K20").Value = 0.72243225772806996
$ws.Range("W21").Value = 0.89691353084840686
$ws.Range("BL21").Value = 0.61805947912757997
$ws.Range("W22").Value = 0.85391851107748884
$ws.Range("AG22").Value = 0.77276300400083553
$ws.Range("BF22").Value = 0.92545164729351059
$ws.Range("BO22").Value = 0.74802842996796248
$ws.Range("AO23").Value = 0.80906432264080674
$ws.Range("AU23").Value = 0.71908111996352142
$ws.Range("AZ23").Value = 0.82328822631810772
$ws.Range("BP23").Value = 0.85691801030710668
$ws.Range("AJ24").Value = 0.54531743015827661
$ws.Range("BE24").Value = 0.68429644861750227
$ws.Range("AI25").Value = 0.71480309712001566
$ws.Range("BB25").Value = 0.81970673566671715
$ws.Range("BF25").Value = 0.89613500032509408
$ws.Range("A26").Value = 0.88168919904526533
$ws.Range("AG26").Value = 0.78825967963918742
$ws.Range("Y27").Value = 0.95796129494590954
$ws.Range("U28").Value = 0.93257637833976048
$ws.Range("AB29").Value = 0.99283767035006654
$ws.Range("BO29").Value = 0.59680920246957636
$ws.Range("E30").Value = 0.66639820956013418
$ws.Range("I31").Value = 0.85088262165380379
$ws.Range("BL31").Value = 0.94073564589696357
$ws.Range("V32").Value = 0.96288459149571559
$ws.Range("AP32").Value = 0.82733204831268914
$ws.Range("R33").Value = 0.88981614037001555
$ws.Range("U34").Value = 0.86517113324668837
$ws.Range("AF34").Value = 0.85135729681025951
$ws.Range("AW34").Value = 0.90637606056650577
$ws.Range("BG34").Value = 0.86763954302225543
$ws.Range("A35").Value = 0.94196899139128099
$ws.Range("P36").Value = 0.99677071804198458
$ws.Range("BN36").Value = 0.93711100860192098
$ws.Range("R37").Value = 0.8610757846318835
$ws.Range("X37").Value = 0.96412795474418544
$ws.Range("AJ37").Value = 0.98556443340881739
$ws.Range("BG37").Value = 0.67006080110773336
$ws.Range("AP38").Value = 0.96265735959168519
$ws.Range("AO39").Value = 0.76023152326660715
$ws.Range("J40").Value = 0.65356820306350871
$ws.Range("S41").Value = 0.75924255959363218
$ws.Range("AF41").Value = 0.67011135805325217
$ws.Range("BL41").Value = 0.59199789896074329
$ws.Range("AZ42").Value = 0.74865098482743586
$ws.Range("O43").Value = 0.59577771355550979
$ws.Range("R43").Value = 0.87912825850374632
$ws.Range("I44").Value = 0.84768506165828306
$ws.Range("T44").Value = 0.99963817815079925
$ws.Range("AR45").Value = 0.77867257389763322
$ws.Range("BE45").Value = 0.70940758167916784
$ws.Range("C46").Value = 0.93113871392950243
$ws.Range("AP46").Value = 0.89037114145116181
$ws.Range("AU46").Value = 0.80869036232059344
$ws.Range("BF46").Value = 0.8813777594929888
$ws.Range("S48").Value = 0.92271409166143048
$ws.Range("AU48").Value = 0.63060423223321749
$ws.Range("AE49").Value = 0.94676318168206675
$ws.Range("AJ49").Value = 0.61427320832481058
$ws.Range("AU49").Value = 0.90706103142396555
$ws.Range("AL50").Value = 0.96874116160289137
$ws.Range("AQ50").Value = 0.88156081985337287
$ws.Range("Q51").Value = 0.61412833171173675
$ws.Range("AD51").Value = 0.62638092938739032
$ws.Range("AJ51").Value = 0.66248809396463826
$ws.Range("AN51").Value = 0.86815704761491252
$ws.Range("AG52").Value = 0.74736357809848408
$ws.Range("BH52").Value = 0.75926510840963335
$ws.Range("H53").Value = 0.89926629707459704
$ws.Range("R53").Value = 0.70168994343982627
$ws.Range("AN53").Value = 0.93437624470190594
$ws.Range("BG53").Value = 0.87625060520719489
$ws.Range("A54").Value = 0.97079998692500791
$ws.Range("AU55").Value = 0.99429476067658218
$ws.Range("BO55").Value = 0.8380864353578944
$ws.Range("T56").Value = 0.76791238373627912
$ws.Range("AM56").Value = 0.82358842522109543
$ws.Range("BF56").Value = 0.77362073465324221
$ws.Range("AN59").Value = 0.75796304141242943
$ws.Range("P60").Value = 0.8021266356174237
$ws.Range("AU62").Value = 0.73087421378009176
$ws.Range("AV62").Value = 0.99091488244356352
$ws.Range("BI62").Value = 0.93267346647763227
$ws.Range("BL62").Value = 0.96149430280741632
$ws.Range("L63").Value = 0.92598967772275176
$ws.Range("BI63").Value = 0.91554353474077277
$ws.Range("BM63").Value = 0.78160487768006404
$ws.Range("B64").Value = 0.976241462136751
$ws.Range("K64").Value = 0.81650424107506425
$ws.Range("B65").Value = 0.85407755698835452
$ws.Range("N65").Value = 0.91773194242302814
$ws.Range("BN65").Value = 0.9512864376028709
$ws.Range("BN68").Value = 0.74445071798592799
